$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13 for the GIRAFE product, shifting rows 13-42 down to 14-43
$ws.Rows.Item(13).Insert()

# Restore column-A formatting (border/bold/center) on the newly inserted row, copied from row below
$ws.Range('A14').Copy()
$ws.Range('A13').PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 2
$ws.Range('A2').Value = 0
$ws.Range('B2').Value = 'Aerosols'
$ws.Range('C2').Value = 'Aerosols'
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '1995-06-01'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '2023-10-31'
$ws.Range('E2').Style = 'Normal'
$ws.Range('F2').Value = 'Atmospheric Composition'

# Row 3
$ws.Range('A3').Value = 2
$ws.Range('B3').Value = 'CH4'
$ws.Range('C3').Value = 'Greenhouse Gases'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2002-10-01'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '2022-12-31'
$ws.Range('E3').Style = 'Normal'
$ws.Range('F3').Value = 'Atmospheric Composition'

# Row 4
$ws.Range('A4').Value = 1
$ws.Range('B4').Value = 'CO2'
$ws.Range('C4').Value = 'Greenhouse Gases'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '2002-10-01'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '2022-12-31'
$ws.Range('E4').Style = 'Normal'
$ws.Range('F4').Value = 'Atmospheric Composition'

# Row 5
$ws.Range('A5').Value = 3
$ws.Range('B5').Value = 'Ozone'
$ws.Range('C5').Value = 'Ozone'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '1970-04-01'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '2023-07-31'
$ws.Range('E5').Style = 'Normal'
$ws.Range('F5').Value = 'Atmospheric Composition'

# Row 6
$ws.Range('A6').Value = 4
$ws.Range('B6').Value = 'Clouds'
$ws.Range('C6').Value = 'Clouds'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1979-01-01'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '2023-10-31'
$ws.Range('E6').Style = 'Normal'
$ws.Range('F6').Value = 'Atmospheric Physics'

# Row 7
$ws.Range('A7').Value = 8
$ws.Range('B7').Value = 'Earth Rad Budget - C3S CCI'
$ws.Range('C7').Value = 'Earth Radiation Budget'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1979-01-01'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '2024-09-17'
$ws.Range('E7').Style = 'Normal'
$ws.Range('F7').Value = 'Atmospheric Physics'

# Row 8
$ws.Range('A8').Value = 9
$ws.Range('B8').Value = 'Earth Rad Budget - C3S RMIB TotSolarIrrad'
$ws.Range('C8').Value = 'Earth Radiation Budget'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1979-01-01'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '2024-12-13'
$ws.Range('E8').Style = 'Normal'
$ws.Range('F8').Value = 'Atmospheric Physics'

# Row 9
$ws.Range('A9').Value = 7
$ws.Range('B9').Value = 'Earth Rad Budget - ESA CLOUD CCI'
$ws.Range('C9').Value = 'Earth Radiation Budget'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1979-01-01'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '2024-09-17'
$ws.Range('E9').Style = 'Normal'
$ws.Range('F9').Value = 'Atmospheric Physics'

# Row 10
$ws.Range('A10').Value = 5
$ws.Range('B10').Value = 'Earth Rad Budget - NASA CERES EBAF'
$ws.Range('C10').Value = 'Earth Radiation Budget'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1979-01-01'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '2024-09-17'
$ws.Range('E10').Style = 'Normal'
$ws.Range('F10').Value = 'Atmospheric Physics'

# Row 11
$ws.Range('A11').Value = 6
$ws.Range('B11').Value = 'Earth Rad Budget - NOAA/NCEI HIRS'
$ws.Range('C11').Value = 'Earth Radiation Budget'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1979-01-01'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '2024-09-17'
$ws.Range('E11').Style = 'Normal'
$ws.Range('F11').Value = 'Atmospheric Physics'

# Row 12
$ws.Range('A12').Value = 10
$ws.Range('B12').Value = 'Precipitation'
$ws.Range('C12').Value = 'Precipitation'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1979-01-01'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '2024-03-31'
$ws.Range('E12').Style = 'Normal'
$ws.Range('F12').Value = 'Atmospheric Physics'

# Row 13
$ws.Range('A13').Value = 12
$ws.Range('B13').Value = 'Precipitation_GIRAFE'
$ws.Range('C13').Value = 'Precipitation'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2002-01-01'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '2022-12-31'
$ws.Range('E13').Style = 'Normal'
$ws.Range('F13').Value = 'Atmospheric Physics'

# Row 14
$ws.Range('A14').Value = 11
$ws.Range('B14').Value = 'Precipitation_microwave'
$ws.Range('C14').Value = 'Precipitation'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2000-01-01'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '2017-12-31'
$ws.Range('E14').Style = 'Normal'
$ws.Range('F14').Value = 'Atmospheric Physics'

# Row 15
$ws.Range('A15').Value = 14
$ws.Range('B15').Value = 'Surface Rad Budget (CCI)'
$ws.Range('C15').Value = 'Surface Radiation Budget'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1979-01-01'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '2023-10-31'
$ws.Range('E15').Style = 'Normal'
$ws.Range('F15').Value = 'Atmospheric Physics'

# Row 16
$ws.Range('A16').Value = 13
$ws.Range('B16').Value = 'Surface Rad Budget - (CMSAF/CLARA)'
$ws.Range('C16').Value = 'Surface Radiation Budget'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1979-01-01'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '2023-10-31'
$ws.Range('E16').Style = 'Normal'
$ws.Range('F16').Value = 'Atmospheric Physics'

# Row 17
$ws.Range('A17').Value = 15
$ws.Range('B17').Value = 'Total Column Water Vapour (HOAPS)'
$ws.Range('C17').Value = 'Upper-air Water Vapour'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1988-01-31'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '2020-12-31'
$ws.Range('E17').Style = 'Normal'
$ws.Range('F17').Value = 'Atmospheric Physics'

# Row 18
$ws.Range('A18').Value = 16
$ws.Range('B18').Value = 'Total Column Water Vapour (MERIS/SSMI)'
$ws.Range('C18').Value = 'Upper-air Water Vapour'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2002-05-01'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '2017-12-31'
$ws.Range('E18').Style = 'Normal'
$ws.Range('F18').Value = 'Atmospheric Physics'

# Row 19
$ws.Range('A19').Value = 17
$ws.Range('B19').Value = 'Tropospheric Humidity Profiles (RO)'
$ws.Range('C19').Value = 'Upper-air Water Vapour'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2006-12-01'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '2024-02-29'
$ws.Range('E19').Style = 'Normal'
$ws.Range('F19').Value = 'Atmospheric Physics'

# Row 20
$ws.Range('A20').Value = 18
$ws.Range('B20').Value = 'Upper Tropospheric Humidity'
$ws.Range('C20').Value = 'Upper-air Water Vapour'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1999-01-01'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '2021-02-28'
$ws.Range('E20').Style = 'Normal'
$ws.Range('F20').Value = 'Atmospheric Physics'

# Row 21
$ws.Range('A21').Value = 23
$ws.Range('B21').Value = 'Glaciers elevation and mass change data'
$ws.Range('C21').Value = 'Glaciers'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1975-04-01'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '2021-09-30'
$ws.Range('E21').Style = 'Normal'
$ws.Range('F21').Value = 'Cryosphere'

# Row 22
$ws.Range('A22').Value = 24
$ws.Range('B22').Value = 'Randolph Glacier Inventory for the year 2000'
$ws.Range('C22').Value = 'Glaciers'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1990-01-01'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '2010-12-31'
$ws.Range('E22').Style = 'Normal'
$ws.Range('F22').Value = 'Cryosphere'

# Row 23
$ws.Range('A23').Value = 20
$ws.Range('B23').Value = 'Ice Sheet Gravimetric Mass Balance'
$ws.Range('C23').Value = 'Ice Sheets'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2002-04-16'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '2022-12-17'
$ws.Range('E23').Style = 'Normal'
$ws.Range('F23').Value = 'Cryosphere'

# Row 24
$ws.Range('A24').Value = 21
$ws.Range('B24').Value = 'Ice Sheet Surface Elevation Change (Antarctica)'
$ws.Range('C24').Value = 'Ice Sheets'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1994-11-01'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '2020-06-01'
$ws.Range('E24').Style = 'Normal'
$ws.Range('F24').Value = 'Cryosphere'

# Row 25
$ws.Range('A25').Value = 22
$ws.Range('B25').Value = 'Ice Sheet Surface Elevation Change (Greenland)'
$ws.Range('C25').Value = 'Ice Sheets'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1992-01-01'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '2024-01-01'
$ws.Range('E25').Style = 'Normal'
$ws.Range('F25').Value = 'Cryosphere'

# Row 26
$ws.Range('A26').Value = 19
$ws.Range('B26').Value = 'Ice Sheet Velocity (Greenland)'
$ws.Range('C26').Value = 'Ice Sheets'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2018-10-01'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '2021-09-30'
$ws.Range('E26').Style = 'Normal'
$ws.Range('F26').Value = 'Cryosphere'

# Row 27
$ws.Range('A27').Value = 33
$ws.Range('B27').Value = 'Surface Albedo 10-daily'
$ws.Range('C27').Value = 'Albedo'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1981-09-20'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '2020-06-30'
$ws.Range('E27').Style = 'Normal'
$ws.Range('F27').Value = 'Land Biosphere'

# Row 28
$ws.Range('A28').Value = 29
$ws.Range('B28').Value = 'FAPAR'
$ws.Range('C28').Value = 'FAPAR'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1981-09-20'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '2020-06-30'
$ws.Range('E28').Style = 'Normal'
$ws.Range('F28').Value = 'Land Biosphere'

# Row 29
$ws.Range('A29').Value = 30
$ws.Range('B29').Value = 'Fire Burned Areas'
$ws.Range('C29').Value = 'Fire'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2001-01-01'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '2022-12-01'
$ws.Range('E29').Style = 'Normal'
$ws.Range('F29').Value = 'Land Biosphere'

# Row 30
$ws.Range('A30').Value = 31
$ws.Range('B30').Value = 'Fire Radiative Power'
$ws.Range('C30').Value = 'Fire'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2020-01-01'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '2024-02-29'
$ws.Range('E30').Style = 'Normal'
$ws.Range('F30').Value = 'Land Biosphere'

# Row 31
$ws.Range('A31').Value = 28
$ws.Range('B31').Value = 'LAI'
$ws.Range('C31').Value = 'LAI'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1981-09-20'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '2020-06-30'
$ws.Range('E31').Style = 'Normal'
$ws.Range('F31').Value = 'Land Biosphere'

# Row 32
$ws.Range('A32').Value = 32
$ws.Range('B32').Value = 'Land Cover'
$ws.Range('C32').Value = 'Land Cover'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1992-01-01'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '2022-12-31'
$ws.Range('E32').Style = 'Normal'
$ws.Range('F32').Value = 'Land Biosphere'

# Row 33
$ws.Range('A33').Value = 25
$ws.Range('B33').Value = 'Lake Surface Temperature'
$ws.Range('C33').Value = 'Lakes'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1995-06-01'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '2023-12-31'
$ws.Range('E33').Style = 'Normal'
$ws.Range('F33').Value = 'Land Hydrology'

# Row 34
$ws.Range('A34').Value = 26
$ws.Range('B34').Value = 'Lake Water Level'
$ws.Range('C34').Value = 'Lakes'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1992-10-13'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '2023-12-24'
$ws.Range('E34').Style = 'Normal'
$ws.Range('F34').Value = 'Land Hydrology'

# Row 35
$ws.Range('A35').Value = 27
$ws.Range('B35').Value = 'Soil Moisture'
$ws.Range('C35').Value = 'Soil Moisture'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1978-11-01'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '2024-11-30'
$ws.Range('E35').Style = 'Normal'
$ws.Range('F35').Value = 'Land Hydrology'

# Row 36
$ws.Range('A36').Value = 34
$ws.Range('B36').Value = 'Ocean Colour'
$ws.Range('C36').Value = 'Ocean Colour'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1997-09-04'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '2024-09-30'
$ws.Range('E36').Style = 'Normal'
$ws.Range('F36').Value = 'Ocean'

# Row 37
$ws.Range('A37').Value = 38
$ws.Range('B37').Value = 'SST'
$ws.Range('C37').Value = 'SST'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1981-08-24'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '2022-12-31'
$ws.Range('E37').Style = 'Normal'
$ws.Range('F37').Value = 'Ocean'

# Row 38
$ws.Range('A38').Value = 39
$ws.Range('B38').Value = 'SST (ESA CCI GMPE)'
$ws.Range('C38').Value = 'SST'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1981-09-01'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '2016-12-31'
$ws.Range('E38').Style = 'Normal'
$ws.Range('F38').Value = 'Ocean'

# Row 39
$ws.Range('A39').Value = 35
$ws.Range('B39').Value = 'Sea Ice Concentration'
$ws.Range('C39').Value = 'Sea Ice'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1978-10-25'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '2024-12-02'
$ws.Range('E39').Style = 'Normal'
$ws.Range('F39').Value = 'Ocean'

# Row 40
$ws.Range('A40').Value = 36
$ws.Range('B40').Value = 'Sea Ice Edge and Type'
$ws.Range('C40').Value = 'Sea Ice'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1978-10-25'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '2024-12-02'
$ws.Range('E40').Style = 'Normal'
$ws.Range('F40').Value = 'Ocean'

# Row 41
$ws.Range('A41').Value = 37
$ws.Range('B41').Value = 'Sea Ice Thickness'
$ws.Range('C41').Value = 'Sea Ice'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2002-10-01'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '2024-04-30'
$ws.Range('E41').Style = 'Normal'
$ws.Range('F41').Value = 'Ocean'

# Row 42
$ws.Range('A42').Value = 40
$ws.Range('B42').Value = 'Sea Level'
$ws.Range('C42').Value = 'Sea Level'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1993-01-01'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '2023-12-31'
$ws.Range('E42').Style = 'Normal'
$ws.Range('F42').Value = 'Ocean'

# Row 43
$ws.Range('A43').Value = 41
$ws.Range('B43').Value = 'Surface Geostrophic Currents'
$ws.Range('C43').Value = 'Surface Currents'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1993-01-01'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '2023-12-31'
$ws.Range('E43').Style = 'Normal'
$ws.Range('F43').Value = 'Ocean'
